$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style/formatting of the existing header row (bold, centered, bordered)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the team record (Wins/Losses/Ties) for every data row (2 through 39)
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 30).Value = 96   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 66   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
